$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (swap H5 and H6 weight values)
$ws.Range("H5").Value = 0.15
$ws.Range("H6").Value = 0.35

# Update column widths
$ws.Columns.Item(1).ColumnWidth = 21.9595141700405
$ws.Columns.Item(4).ColumnWidth = 78.0890688259109

# Update selection / top-left cell (view state)
$ws.Range("H7").Select()
$excel.ActiveWindow.ScrollColumn = 4

# Update tab ratio (workbook window zoom ratio for sheet tabs)
$excel.ActiveWindow.TabRatio = 989
